# Scheduled market-data refresh: update leve profit calculations (currentAveragePrice*,
# LevePriceNQ/HQ, LeveProfitNQ/HQ) across the ALC, ARM, CRP, CUL, GSM, LTW and WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1097296.9
$ws.Range("I19").Value = 1754963.1
$ws.Range("J19").Value = 1186.6666
$ws.Range("K19").Value = 1754963.1
$ws.Range("L19").Value = 1186.6666
$ws.Range("M19").Value = -1754788.1
$ws.Range("N19").Value = -1536.6666
$ws.Range("H33").Value = 347.1
$ws.Range("I33").Value = 296.125
$ws.Range("K33").Value = 296.125
$ws.Range("M33").Value = -67.125
$ws.Range("H101").Value = 1721.5
$ws.Range("J101").Value = 1964.2222
$ws.Range("L101").Value = 5892.6666
$ws.Range("N101").Value = -9136.6666
$ws.Range("H137").Value = 629489.75
$ws.Range("I137").Value = 1987836.8
$ws.Range("J137").Value = 2560.3845
$ws.Range("K137").Value = 5963510.4
$ws.Range("L137").Value = 7681.1535
$ws.Range("M137").Value = -5960960.4
$ws.Range("N137").Value = -12781.1535
$ws.Range("H138").Value = 2984.6667
$ws.Range("I138").Value = 1481.7084
$ws.Range("J138").Value = 3909.5642
$ws.Range("K138").Value = 4445.1252
$ws.Range("L138").Value = 11728.6926
$ws.Range("M138").Value = 694.8747999999996
$ws.Range("N138").Value = -22008.6926

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H94").Value = 34933.332
$ws.Range("J94").Value = 34933.332
$ws.Range("L94").Value = 34933.332
$ws.Range("N94").Value = -36735.332
$ws.Range("H104").Value = 33500
$ws.Range("J104").Value = 33500
$ws.Range("L104").Value = 33500
$ws.Range("N104").Value = -40488
$ws.Range("H122").Value = 3249.2
$ws.Range("I122").Value = 2944.1614
$ws.Range("J122").Value = 4299.8887
$ws.Range("K122").Value = 8832.484199999999
$ws.Range("L122").Value = 12899.6661
$ws.Range("M122").Value = -6382.484199999999
$ws.Range("N122").Value = -17799.6661
$ws.Range("H132").Value = 4595.8
$ws.Range("I132").Value = 2836.6667
$ws.Range("K132").Value = 8510.000100000001
$ws.Range("M132").Value = -5980.000100000001
$ws.Range("H137").Value = 42304.832
$ws.Range("J137").Value = 42304.832
$ws.Range("L137").Value = 42304.832
$ws.Range("N137").Value = -52504.832

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1895.9727
$ws.Range("I31").Value = 863.3333
$ws.Range("J31").Value = 2616.4187
$ws.Range("K31").Value = 863.3333
$ws.Range("L31").Value = 2616.4187
$ws.Range("M31").Value = -568.3333
$ws.Range("N31").Value = -3206.4187
$ws.Range("H34").Value = 1895.9727
$ws.Range("I34").Value = 863.3333
$ws.Range("J34").Value = 2616.4187
$ws.Range("K34").Value = 863.3333
$ws.Range("L34").Value = 2616.4187
$ws.Range("M34").Value = -661.3333
$ws.Range("N34").Value = -3020.4187
$ws.Range("H64").Value = 26517.5
$ws.Range("J64").Value = 26517.5
$ws.Range("L64").Value = 26517.5
$ws.Range("N64").Value = -27013.5
$ws.Range("H67").Value = 26517.5
$ws.Range("J67").Value = 26517.5
$ws.Range("L67").Value = 26517.5
$ws.Range("N67").Value = -28233.5
$ws.Range("H86").Value = 2822.4119
$ws.Range("I86").Value = 2441.5715
$ws.Range("K86").Value = 2441.5715
$ws.Range("M86").Value = -1318.5715
$ws.Range("H89").Value = 2822.4119
$ws.Range("I89").Value = 2441.5715
$ws.Range("K89").Value = 12207.8575
$ws.Range("M89").Value = -6591.8575

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 496082.78
$ws.Range("I5").Value = 1348.3636
$ws.Range("J5").Value = 836212.7
$ws.Range("K5").Value = 4045.0908
$ws.Range("L5").Value = 2508638.1
$ws.Range("M5").Value = -3933.0908
$ws.Range("N5").Value = -2508862.1
$ws.Range("H68").Value = 1625
$ws.Range("I68").Value = 1409.2727
$ws.Range("J68").Value = 1684.325
$ws.Range("K68").Value = 4227.8181
$ws.Range("L68").Value = 5052.975
$ws.Range("M68").Value = -3416.8181
$ws.Range("N68").Value = -6674.975
$ws.Range("H71").Value = 1625
$ws.Range("I71").Value = 1409.2727
$ws.Range("J71").Value = 1684.325
$ws.Range("K71").Value = 12683.4543
$ws.Range("L71").Value = 15158.925
$ws.Range("M71").Value = -8627.454299999999
$ws.Range("N71").Value = -23270.925
$ws.Range("H113").Value = 5000610
$ws.Range("I113").Value = 643.61536
$ws.Range("J113").Value = 10417240
$ws.Range("K113").Value = 1930.84608
$ws.Range("L113").Value = 31251720
$ws.Range("M113").Value = 239.15392
$ws.Range("N113").Value = -31256060
$ws.Range("H122").Value = 2227.6482
$ws.Range("I122").Value = 719.9
$ws.Range("J122").Value = 3114.5588
$ws.Range("K122").Value = 6479.099999999999
$ws.Range("L122").Value = 28031.0292
$ws.Range("M122").Value = -4029.099999999999
$ws.Range("N122").Value = -32931.0292
$ws.Range("H131").Value = 804.8889
$ws.Range("I131").Value = 362
$ws.Range("J131").Value = 828.44684
$ws.Range("K131").Value = 1086
$ws.Range("L131").Value = 2485.34052
$ws.Range("M131").Value = 3954
$ws.Range("N131").Value = -12565.34052
$ws.Range("H135").Value = 496082.78
$ws.Range("I135").Value = 1348.3636
$ws.Range("J135").Value = 836212.7
$ws.Range("K135").Value = 12135.2724
$ws.Range("L135").Value = 7525914.3
$ws.Range("M135").Value = -9600.2724
$ws.Range("N135").Value = -7530984.3
$ws.Range("H139").Value = 802.1053000000001
$ws.Range("I139").Value = 802.1053000000001
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 2406.3159
$ws.Range("L139").Value = 0
$ws.Range("M139").Value = 2733.6841
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H88").Value = 32388.234
$ws.Range("J88").Value = 32388.234
$ws.Range("L88").Value = 32388.234
$ws.Range("N88").Value = -33290.234
$ws.Range("H91").Value = 32388.234
$ws.Range("J91").Value = 32388.234
$ws.Range("L91").Value = 32388.234
$ws.Range("N91").Value = -35508.234
$ws.Range("H97").Value = 1339.2222
$ws.Range("I97").Value = 1382
$ws.Range("J97").Value = 1285.75
$ws.Range("K97").Value = 1382
$ws.Range("L97").Value = 1285.75
$ws.Range("M97").Value = -886
$ws.Range("N97").Value = -2277.75
$ws.Range("H102").Value = 2723.1875
$ws.Range("I102").Value = 1943
$ws.Range("J102").Value = 5063.75
$ws.Range("K102").Value = 1943
$ws.Range("L102").Value = 5063.75
$ws.Range("M102").Value = -321
$ws.Range("N102").Value = -8307.75
$ws.Range("H123").Value = 10560
$ws.Range("J123").Value = 10560
$ws.Range("L123").Value = 10560
$ws.Range("N123").Value = -15460

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H14").Value = 28031.8
$ws.Range("I14").Value = 8000
$ws.Range("J14").Value = 30257.555
$ws.Range("K14").Value = 8000
$ws.Range("L14").Value = 30257.555
$ws.Range("M14").Value = -7828
$ws.Range("N14").Value = -30601.555
$ws.Range("H132").Value = 4271.9062
$ws.Range("I132").Value = 3442.0952
$ws.Range("J132").Value = 5856.091
$ws.Range("K132").Value = 10326.2856
$ws.Range("L132").Value = 17568.273
$ws.Range("M132").Value = -7796.285600000001
$ws.Range("N132").Value = -22628.273
$ws.Range("H133").Value = 31858.75
$ws.Range("J133").Value = 31858.75
$ws.Range("L133").Value = 31858.75
$ws.Range("N133").Value = -36918.75
$ws.Range("H134").Value = 39602.133
$ws.Range("J134").Value = 39602.133
$ws.Range("L134").Value = 39602.133
$ws.Range("N134").Value = -49742.133
$ws.Range("H135").Value = 94000
$ws.Range("J135").Value = 94000
$ws.Range("L135").Value = 94000
$ws.Range("N135").Value = -104140

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H104").Value = 30000
$ws.Range("J104").Value = 30000
$ws.Range("L104").Value = 30000
$ws.Range("N104").Value = -36988
$ws.Range("H107").Value = 1330.8334
$ws.Range("I107").Value = 592.5
$ws.Range("J107").Value = 1700
$ws.Range("K107").Value = 1777.5
$ws.Range("L107").Value = 5100
$ws.Range("M107").Value = 142.5
$ws.Range("N107").Value = -8940
$ws.Range("H113").Value = 308.21054
$ws.Range("I113").Value = 306.25
$ws.Range("J113").Value = 311.57144
$ws.Range("K113").Value = 918.75
$ws.Range("L113").Value = 934.71432
$ws.Range("M113").Value = 1251.25
$ws.Range("N113").Value = -5274.71432
$ws.Range("H122").Value = 3932.5356
$ws.Range("I122").Value = 2554.2222
$ws.Range("J122").Value = 4585.421
$ws.Range("K122").Value = 7662.6666
$ws.Range("L122").Value = 13756.263
$ws.Range("M122").Value = -5212.6666
$ws.Range("N122").Value = -18656.263
